# Update OleDb connection string
# On the "LoginUser" sheet, insert a new "LoginAsAdmin" test row right after the
# header row, and remove the old "SuccessfulLogIn" row at the bottom.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LoginUser")

# Remove the last data row (row 5: SuccessfulLogIn / georgi_vatashki@abv.bg / softuni123)
$ws.Rows.Item(5).Delete()

# Insert a new row above the current row 2 (LoginWithoutEmail), shifting data down
$ws.Rows.Item(2).Insert()

# Reset the freshly inserted row to the default (general) style before filling it in
# so the numeric cell below is stored as a real number rather than inheriting the
# header row's text format.
$ws.Range("A2:C2").Style = "Normal"

# Fill in the new row 2 with the LoginAsAdmin test case (email before name so the
# shared-string table reuses slots in the same order as the source workbook).
$ws.Cells.Item(2, 2).Value = "admin@admin.com"
$ws.Cells.Item(2, 1).Value = "LoginAsAdmin"
$ws.Cells.Item(2, 3).Value = 123

# Match formatting of neighboring data rows: TestName/Number columns are centered.
$ws.Cells.Item(2, 1).HorizontalAlignment = -4108
$ws.Cells.Item(2, 3).HorizontalAlignment = -4108

# Update selected cell as recorded in the workbook view
$ws.Range("A2").Select()
